$d = $word.ActiveDocument
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Paragraph "Pruebas de la aplicación." -> add <w:lastRenderedPageBreak/>
#    before the <w:t> in its run.
# ---------------------------------------------------------------------------
$pPruebas = $d.Paragraphs.Item(43)
$rPruebas = $pPruebas.Range
$xmlPruebas = '<w:p xmlns:w="' + $w + '"><w:r><w:lastRenderedPageBreak/><w:t>Pruebas de la aplicaci' + [char]0x00F3 + 'n.</w:t></w:r></w:p>'
$rPruebas.InsertXML($xmlPruebas) | Out-Null

# ---------------------------------------------------------------------------
# 2) Paragraph "Construcción de la aplicación." -> new text, and gains the
#    _GoBack bookmark that used to sit at the end of the "Palabra Clave"
#    paragraph.
# ---------------------------------------------------------------------------
$pConstr = $d.Paragraphs.Item(40)
$rConstr = $pConstr.Range
$newConstrText = "En esta fase se construira la aplicaci" + [char]0x00F3 + "n web que permita realizar las busquedas inteligentes de sitios de interes y direcciones urbanas, para ello sera desarollada en JEE y se utlizara la librer" + [char]0x00ED + "a JENA para hacer consultas con el lenguaje SPARQL a la ontologia creada"
$xmlConstr = '<w:p xmlns:w="' + $w + '"><w:r><w:t>' + $newConstrText + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rConstr.InsertXML($xmlConstr) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the old _GoBack bookmark (previously at the end of the
#    "Palabra Clave" paragraph).
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4) Delete the empty bold paragraph that used to sit between the
#    "Palabra Clave" section and "2.5. Fase V".
# ---------------------------------------------------------------------------
$pBoldEmpty = $d.Paragraphs.Item(38)
$pBoldEmpty.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 5) "Clase Palabra Clave" paragraph -> split the trailing run so that
#    "sinonimo" becomes "palabra" as its own run, followed by a separate
#    "." run.
# ---------------------------------------------------------------------------
$pPalabra = $d.Paragraphs.Item(36)
$rPalabra = $pPalabra.Range
$xmlPalabra = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Clase </w:t></w:r>' +
    '<w:r><w:t>Palabra Clave:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Esta clase tiene instancias para palabras clave que describen</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> al tipo de vias definidas y los sitios de interes</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, de tal forma que facilite realizar las busquedas. La caracteristica que define a esta clase es: </w:t></w:r>' +
    '<w:r><w:t>palabra</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
$rPalabra.InsertXML($xmlPalabra) | Out-Null

# ---------------------------------------------------------------------------
# 6) The two blank spacer paragraphs right after "Clase Manzana" and right
#    after "Clase Sitio_Interes" pick up the "Prrafodelista" (List
#    Paragraph) style.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(33).set_Style("List Paragraph")
$d.Paragraphs.Item(35).set_Style("List Paragraph")
